$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "Travel checklist: here's what you need" "Чек-лист путешественника: что вам понадобится"

Replace-Text "Here’s a checklist of the necessary items for your trip: " "Вот контрольный список для вашей поездки: "

Replace-Text "Passport " "Паспорт "

Replace-Text "For travellers from yellow fever endemic countries, follow the requirements set by your country. Vaccination should be done no less than 14 days prior to the journey. " "For travellers from yellow fever endemic countries, follow the requirements set by your country. Вакцинация должна быть проведена не менее чем за 14 дней до поездки. "

Replace-Text "A digital or printed copy of the travel itinerary" "Цифровая или печатная копия маршрута путешествия"

Replace-Text "Smart casual attire for the conference" "Повседневная одежда в стиле «smart-casual» для конференции"

Replace-Text "Black tie attire for the Gala dinner" "Одежда в стиле «черный галстук» для гала-ужина"
